# Apply the "Add files via upload" edit to the MEAM converter workbook.
# This replaces the Pt (platinum) reference-data row with Y (yttrium) data,
# switches the crystal structure from fcc to hcp, and updates the related
# geometry constant (E7/J7) plus the scratch note in K7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: main per-element parameter row -------------------------------
$ws.Range("A3").Value = "Y"
$ws.Range("B3").Value = 4.41
$ws.Range("C3").Value = 3.6

# D3 used to be a formula (=H5); it is now a hard-coded literal value.
$ws.Range("D3").Value = 4.17

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.5

# G3:I3 change value AND number format (from 0.0 to 0.00, matching F3's style)
$ws.Range("G3").NumberFormat = "0.00"
$ws.Range("G3").Value = 0
$ws.Range("H3").NumberFormat = "0.00"
$ws.Range("H3").Value = 0.54
$ws.Range("I3").NumberFormat = "0.00"
$ws.Range("I3").Value = 5

$ws.Range("K3").Value = 6.6
$ws.Range("L3").Value = 14.6
$ws.Range("M3").Value = -10
$ws.Range("N3").Value = 2
$ws.Range("P3").Value = 4
$ws.Range("Q3").Value = 0.1
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0

# --- Row 7: crystal-structure / lattice-geometry row ----------------------
$ws.Range("C7").Value = "hcp"

# E7 used to be a formula (=SQRT(2)); it is now a hard-coded literal value.
$ws.Range("E7").Value = 1

$ws.Range("F7").Value = 2

# J7 used to be a hard-coded literal (1.6); it is now a formula.
$ws.Range("J7").Formula = "=SQRT(8/3)"

# K7 gets an explanatory note next to the new J7 formula.
$ws.Range("K7").Value = "<- sqrt(8/3)"

# --- Row 16: propagate the new number format to the dependent formulas ----
# (C16=G3, D16=H3, E16=I3 mirror B16=F3's "0.00" / no-border style)
$ws.Range("C16").NumberFormat = "0.00"
$ws.Range("D16").NumberFormat = "0.00"
$ws.Range("E16").NumberFormat = "0.00"

# --- Selection / view state, matching the saved worksheet view ------------
$ws.Range("A1").Select()
$ws.Range("L7").Select()
